# Commit: Delete rows [1, 2] from Card24
#
# "Card24" is the data table's 0-indexed data rows 1 and 2 (the header is
# row 1 of the sheet, so data row 0 = sheet row 2, data row 1 = sheet row 3,
# data row 2 = sheet row 4). Deleting them shifts every subsequent row up
# by two, shrinking the used range from A1:L15 to A1:L13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("A3:A4").EntireRow.Delete()
